$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.575.64'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.346.18'
$ws.Range("E3").Value = '  -2.71%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.09'
$ws.Range("E5").Value = '  -3.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.46'
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  -7.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.421'
$ws.Range("E8").Value = '  -6.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.994'
$ws.Range("E10").Value = '  -8.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.342.99'
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.94'
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '97.292.13'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.10'
$ws.Range("E15").Value = '  -5.30%  '
$ws.Range("E16").Value = '  -7.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.975.46'
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.56'
$ws.Range("E18").Value = '  -8.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.336.09'
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.57'
$ws.Range("E20").Value = '  -2.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.520'
$ws.Range("E21").Value = '  -7.38%  '
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '505.84'
$ws.Range("E23").Value = '  -3.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.34'
$ws.Range("E24").Value = '  -5.57%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000199'
$ws.Range("E25").Value = '  -5.41%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.88'
$ws.Range("E26").Value = '  +6.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.78'
$ws.Range("E27").Value = '  -6.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.13'
$ws.Range("E28").Value = '  -8.47%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.522.26'
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.29'
$ws.Range("E30").Value = '  -5.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.142'
$ws.Range("E31").Value = '  -10.68%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.186'
$ws.Range("E33").Value = '  -8.82%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.55'
$ws.Range("E34").Value = '  +8.19%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.554'
$ws.Range("E36").Value = '  -6.97%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.33'
$ws.Range("E37").Value = '  -6.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.00'
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '531.56'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.151'
$ws.Range("E42").Value = '  -3.44%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.41'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.840'
$ws.Range("E44").Value = '  -4.82%  '
$ws.Range("E45").Value = '  +3.25%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.28'
$ws.Range("E46").Value = '  +7.14%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0423'
$ws.Range("E47").Value = '  -4.18%  '
$ws.Range("B48").Value = 'MantraDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.60'
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.51'
$ws.Range("E49").Value = '  -6.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.63'
$ws.Range("E50").Value = '  +6.21%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.37'
$ws.Range("E51").Value = '  -9.98%  '
